$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the frozen-pane split that was anchored at B3 (xSplit=1, ySplit=2)
$win = $excel.ActiveWindow
$win.Split = $false

# The title cell A1:I1 was merged ("Sheet with your disqualified leads"); unmerge it
# before we delete the row so no merged-range artifacts remain.
$ws.Range("A1:I1").UnMerge()

# Drop the title row entirely - this shifts the old header row ("URL"/"Description",
# bold+fill styling) up to become row 1, and every data row shifts up by one too,
# carrying its formatting with it.
$ws.Rows(1).Delete()

# Trim the empty tail rows so only 11 rows (1 header + 10 blank data rows) remain.
$ws.Rows("12:23").Delete()

# Slightly widen the data columns (A:I) to match the refreshed template.
$ws.Range("A1:I11").ColumnWidth = 8.43
